$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.543.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.583.83'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.00%  '
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.578'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.589.45'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.347'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.14%  '
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.041.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.511.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.588.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '345.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.32%  '
$ws.Range('E21').Value = '  +1.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '59.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.420'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0844'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.34'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.73%  '
$ws.Range('E35').Value = '  +2.13%  '
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.857'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.46%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('E41').Value = '  +2.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '295.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.619'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0995'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.997'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0557'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.85'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.999.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.08%  '
